# Fix the "false" boolean value case: the "isSigned" (false) row was stuck
# at the very bottom of the table (row 22), disconnected from its sibling
# boolean field "isCertified" (row 14). Move it so it immediately follows
# "isCertified", and shift the intervening rows (intructionDateTime and the
# funds.* rows) down by one to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: isCertified / isCertified / boolean / TRUE -------------------
# Field/value stay the same; only the "Type" cell now carries the
# highlighted style (s3) that used to live further down the sheet.
$ws.Cells.Item(14,1).Value = "isCertified"
$ws.Cells.Item(14,2).Value = "isCertified"
$ws.Cells.Item(14,3).Value = "boolean"
$ws.Cells.Item(14,4).Value = $true

# --- Row 15: isSigned / isSigned / boolean / FALSE -------------------------
# This used to be row 22 - move it right under isCertified.
$ws.Cells.Item(15,1).Value = "isSigned"
$ws.Cells.Item(15,2).Value = "isSigned"
$ws.Cells.Item(15,3).Value = "boolean"
$ws.Cells.Item(15,4).Value = $false
$ws.Rows.Item(15).RowHeight = 16

# --- Row 16: intructionDateTime / intructionDateTime / str / formula ------
$ws.Cells.Item(16,1).Value = "intructionDateTime"
$ws.Cells.Item(16,2).Value = "intructionDateTime"
$ws.Cells.Item(16,3).Value = "str"
$ws.Cells.Item(16,4).Formula = '=TEXT(NOW(), "yyyy-mm-dd hh:mm:ss.000Z")'
$ws.Rows.Item(16).RowHeight = 17.25

# --- Row 17: funds.1.code / funds.1.code / list.obj.str / ABCD ------------
$ws.Cells.Item(17,1).Value = "funds.1.code"
$ws.Cells.Item(17,2).Value = "funds.1.code"
$ws.Cells.Item(17,3).Value = "list.obj.str"
$ws.Cells.Item(17,4).Value = "ABCD"
$ws.Rows.Item(17).RowHeight = 17.25

# --- Row 18: funds.1.investmentAmt / .. / list.obj.num / 20000 ------------
$ws.Cells.Item(18,1).Value = "funds.1.investmentAmt"
$ws.Cells.Item(18,2).Value = "funds.1.investmentAmt"
$ws.Cells.Item(18,3).Value = "list.obj.num"
$ws.Cells.Item(18,4).Value = 20000
$ws.Rows.Item(18).RowHeight = 17.25

# --- Row 19: funds.1.investmentcury / .. / list.obj.str / HKD -------------
$ws.Cells.Item(19,1).Value = "funds.1.investmentcury"
$ws.Cells.Item(19,2).Value = "funds.1.investmentcury"
$ws.Cells.Item(19,3).Value = "list.obj.str"
$ws.Cells.Item(19,4).Value = "HKD"
$ws.Rows.Item(19).RowHeight = 19.5

# --- Row 20: funds.2.code / funds.2.code / list.obj.str / EFGH ------------
$ws.Cells.Item(20,1).Value = "funds.2.code"
$ws.Cells.Item(20,2).Value = "funds.2.code"
$ws.Cells.Item(20,3).Value = "list.obj.str"
$ws.Cells.Item(20,4).Value = "EFGH"
$ws.Rows.Item(20).RowHeight = 19.5

# --- Row 21: funds.2.investmentAmt / .. / list.obj.num / 1000 -------------
$ws.Cells.Item(21,1).Value = "funds.2.investmentAmt"
$ws.Cells.Item(21,2).Value = "funds.2.investmentAmt"
$ws.Cells.Item(21,3).Value = "list.obj.num"
$ws.Cells.Item(21,4).Value = 1000
$ws.Rows.Item(21).RowHeight = 17.25

# --- Row 22: funds.2.investmentcury / .. / list.obj.str / USD -------------
$ws.Cells.Item(22,1).Value = "funds.2.investmentcury"
$ws.Cells.Item(22,2).Value = "funds.2.investmentcury"
$ws.Cells.Item(22,3).Value = "list.obj.str"
$ws.Cells.Item(22,4).Value = "USD"
$ws.Rows.Item(22).RowHeight = 19.5

# --- Styling ----------------------------------------------------------
# Column C ("Type"): isCertified, isSigned, funds.1.investmentcury and
# funds.2.investmentcury are highlighted (style s3); everything else is
# plain. Copy formats around (instead of stamping raw style indices) so
# the existing style-table entries get reused and no duplicates appear.
$highlightedC = $ws.Cells.Item(22,3)   # already has the highlighted style
$plainC = $ws.Cells.Item(16,3)         # has the default/plain style

$highlightedC.Copy() | Out-Null
$ws.Cells.Item(14,3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15,3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(19,3).PasteSpecial(-4122) | Out-Null

$plainC.Copy() | Out-Null
$ws.Cells.Item(17,3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(18,3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20,3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(21,3).PasteSpecial(-4122) | Out-Null

# Column D ("Value"): numeric amount cells (funds.*.investmentAmt) use the
# right-aligned number style (s2); every other Value cell in this block
# uses the plain style (s1).
$numD = $ws.Cells.Item(3,4)    # amount row, has the numeric style (s2)
$plainD = $ws.Cells.Item(1,4)  # has the plain style (s1)

$plainD.Copy() | Out-Null
$ws.Cells.Item(14,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(16,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(19,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22,4).PasteSpecial(-4122) | Out-Null

$numD.Copy() | Out-Null
$ws.Cells.Item(18,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(21,4).PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# Selection ends up on C14 (as left by the author after the edit).
$ws.Range("C14").Select() | Out-Null
